$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.092565536499023
$ws.Range("B1").Value = 1.824850797653198
$ws.Range("C1").Value = 5.376095294952393
$ws.Range("D1").Value = 0.7536296248435974
$ws.Range("E1").Value = 0.6149309277534485
